# "change abbreviation in company data"
#
# The 항목설명 (item-description) sheet is renumbered/relabelled:
#   - pp / pp_cr / pp_con_num          -> retail / retail_cr / retail_days
#   - fp / fp_cr / fp_con_num          -> foreigner / foreigner_cr / foreigner_days
#   - three brand-new rows are inserted for 기관(Institutional) purchases
#   - itp*  -> invtrust*   (shifted down 3 rows)
#   - penp* -> pension*    (shifted down 3 rows)
#   - pep*  -> privequity* (shifted down 3 rows)
#   - the remaining rows (rc1_pcr .. low_time) shift down 3 rows unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new "Institutional" block: insert 3 rows at 9 ---
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# --- 2. Rename Personal (pp) -> Retail (retail) ---
$ws.Range("B3").Value = "retail"
$ws.Range("C3").Value = "Retail_Purchase_amount (Retail investor, Indivisual investor)"

$ws.Range("B4").Value = "retail_cr"
$ws.Range("C4").Value = "Retail_Purchase_amount_change_rate"

$ws.Range("B5").Value = "retail_days"
$ws.Range("C5").Value = "Retail_Purchase_number_of_consecutive_days"

# --- 3. Rename Foreigner (fp) -> foreigner (abbreviation only) ---
$ws.Range("B6").Value = "foreigner"
$ws.Range("B7").Value = "foreigner_cr"
$ws.Range("B8").Value = "foreigner_days"

# --- 4. Fill the newly inserted Institutional rows ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "institution"
$ws.Range("C9").Value = "Institutional_Purchase_amount"
$ws.Range("D9").Value = "기관매수액"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "institution_cr"
$ws.Range("C10").Value = "Institutional_Purchase_amount_change_rate"
$ws.Range("D10").Value = "기관매수변동률"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "institution_days"
$ws.Range("C11").Value = "Institutional_Purchase_amount_number_of_consecutive_days"
$ws.Range("D11").Value = "기관매수지속일"

# --- 5. Renumber + rename the rows that got pushed down by the insert ---
# InvestmentTrust (itp -> invtrust), now rows 12-14
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "invtrust"
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "invtrust_cr"
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "invtrust_days"

# Pension (penp -> pension), now rows 15-17
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "pension"
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "pension_cr"
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "pension_days"

# PrivateEquity (pep -> privequity), now rows 18-20
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "privequity"
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "privequity_cr"
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "privequity_days"

# --- 6. Renumber the trailing, otherwise-unchanged rows (rc1_pcr .. low_time) ---
# The "번호" numbering restarts at 17 here (matches source data), rows 21-33
$ws.Range("A21").Value = 17
$ws.Range("A22").Value = 18
$ws.Range("A23").Value = 19
$ws.Range("A24").Value = 20
$ws.Range("A25").Value = 21
$ws.Range("A26").Value = 22
$ws.Range("A27").Value = 23
$ws.Range("A28").Value = 24
$ws.Range("A29").Value = 25
$ws.Range("A30").Value = 26
$ws.Range("A31").Value = 27
$ws.Range("A32").Value = 28
$ws.Range("A33").Value = 29

# --- 7. Restore the view state recorded in the saved workbook ---
$ws.Range("C29").Select()
